# Update Name of Algo
# Apply numeric corrections to the RandomForest imputation result sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.212599999999994
$ws.Range("A9").Value = -20.44679999999997
$ws.Range("D11").Value = -8.356100000000003
$ws.Range("A18").Value = -23.07310000000001
$ws.Range("A20").Value = -22.15570000000003
$ws.Range("E21").Value = 12.99109999999999
